$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5526315789473685
$ws.Range("B3").Value = 0.5783582089552238
$ws.Range("B4").Value = 0.4769230769230769
$ws.Range("B5").Value = 0.4057591623036649
$ws.Range("B6").Value = 505.6701173782349
